# Auto-applies updated coin price / 1h-volume figures scraped on 2023-02-16.
# Source cells store plain text (e.g. "320.93", "8.11%") rather than numbers,
# so we force a Text number format before writing the value and then restore
# the "Normal" style so no stray formatting is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "320.93"
$cell.Style = "Normal"
$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = "8.11%"
$cell.Style = "Normal"
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "51.53"
$cell.Style = "Normal"
$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = "23.06%"
$cell.Style = "Normal"
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "5.360"
$cell.Style = "Normal"
$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = "7.12%"
$cell.Style = "Normal"
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "0.08112"
$cell.Style = "Normal"
$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = "7.86%"
$cell.Style = "Normal"
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "4.583"
$cell.Style = "Normal"
$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = "4.60%"
$cell.Style = "Normal"
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "1.658"
$cell.Style = "Normal"
$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = "4.93%"
$cell.Style = "Normal"
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "1.128"
$cell.Style = "Normal"
$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = "21.63%"
$cell.Style = "Normal"
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.1316"
$cell.Style = "Normal"
$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = "10.18%"
$cell.Style = "Normal"
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.1960"
$cell.Style = "Normal"
$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = "7.10%"
$cell.Style = "Normal"
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.09614"
$cell.Style = "Normal"
$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = "7.92%"
$cell.Style = "Normal"
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.04576"
$cell.Style = "Normal"
$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = "11.95%"
$cell.Style = "Normal"
$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = "-0.03%"
$cell.Style = "Normal"
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "0.001323"
$cell.Style = "Normal"
$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = "2.87%"
$cell.Style = "Normal"
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.005788"
$cell.Style = "Normal"
$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = "-3.04%"
$cell.Style = "Normal"
$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = "0.55%"
$cell.Style = "Normal"
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "2.427"
$cell.Style = "Normal"
$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = "1.09%"
$cell.Style = "Normal"
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "8.192"
$cell.Style = "Normal"
$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = "1.53%"
$cell.Style = "Normal"
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "0.1399"
$cell.Style = "Normal"
$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = "0.52%"
$cell.Style = "Normal"
$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = "-11.39%"
$cell.Style = "Normal"
$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = "4.63%"
$cell.Style = "Normal"
$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = "3.13%"
$cell.Style = "Normal"
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "0.004296"
$cell.Style = "Normal"
$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = "10.50%"
$cell.Style = "Normal"
$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = "9.42%"
$cell.Style = "Normal"
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "0.0003717"
$cell.Style = "Normal"
$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = "-0.18%"
$cell.Style = "Normal"
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.02765"
$cell.Style = "Normal"
$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = "14.76%"
$cell.Style = "Normal"
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.05549"
$cell.Style = "Normal"
$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = "6.68%"
$cell.Style = "Normal"
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.006193"
$cell.Style = "Normal"
$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = "-1.83%"
$cell.Style = "Normal"
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.007774"
$cell.Style = "Normal"
$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = "-0.36%"
$cell.Style = "Normal"
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.1445"
$cell.Style = "Normal"
$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = "9.07%"
$cell.Style = "Normal"
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.007679"
$cell.Style = "Normal"
$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = "3.75%"
$cell.Style = "Normal"
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.008828"
$cell.Style = "Normal"
$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = "18.21%"
$cell.Style = "Normal"
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.3525"
$cell.Style = "Normal"
$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = "19.34%"
$cell.Style = "Normal"
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.00006815"
$cell.Style = "Normal"
$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = "5.87%"
$cell.Style = "Normal"
$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = "-0.25%"
$cell.Style = "Normal"
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.06125"
$cell.Style = "Normal"
$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = "82.28%"
$cell.Style = "Normal"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.003995"
$cell.Style = "Normal"
$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = "-5.00%"
$cell.Style = "Normal"
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.00002098"
$cell.Style = "Normal"
$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = "-0.25%"
$cell.Style = "Normal"
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.0001998"
$cell.Style = "Normal"
$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = "-0.25%"
$cell.Style = "Normal"
